# Weekly Fruit/Vegetable price update: insert two new rows of data
# (Vega Monumental Concepción - Palta) ahead of the existing history,
# pushing the prior rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 267/268; everything currently at row 267 and
# below shifts down by two rows (267->269, ..., 295->297).
$ws.Rows.Item(267).Insert()
$ws.Rows.Item(268).Insert()

# --- New row 267: Hass / "1a nueva(o)" / Cabildo ---
$ws.Range("A267").Value = 11
$ws.Range("B267").Value = "Vega Monumental Concepción"
$ws.Range("C267").Value = "Bíobío"
$ws.Range("D267").Value = 44461
$ws.Range("E267").Value = 8
$ws.Range("F267").Value = "Fruta"
$ws.Range("G267").Value = 100106
$ws.Range("H267").Value = "Oleaginosos"
$ws.Range("I267").Value = 100106002
$ws.Range("J267").Value = "Palta"
$ws.Range("K267").Value = "Hass"
$ws.Range("L267").Value = "1a nueva(o)"
$ws.Range("M267").Value = 50
$ws.Range("N267").Value = 3000
$ws.Range("O267").Value = 3000
$ws.Range("P267").Value = 3000
$ws.Range("Q267").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R267").Value = "Cabildo"
$ws.Range("S267").Value = 3000
$ws.Range("T267").Value = 1

# --- New row 268: Hass / "2a nueva(o)" / Cabildo ---
$ws.Range("A268").Value = 11
$ws.Range("B268").Value = "Vega Monumental Concepción"
$ws.Range("C268").Value = "Bíobío"
$ws.Range("D268").Value = 44461
$ws.Range("E268").Value = 8
$ws.Range("F268").Value = "Fruta"
$ws.Range("G268").Value = 100106
$ws.Range("H268").Value = "Oleaginosos"
$ws.Range("I268").Value = 100106002
$ws.Range("J268").Value = "Palta"
$ws.Range("K268").Value = "Hass"
$ws.Range("L268").Value = "2a nueva(o)"
$ws.Range("M268").Value = 50
$ws.Range("N268").Value = 2800
$ws.Range("O268").Value = 2800
$ws.Range("P268").Value = 2800
$ws.Range("Q268").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R268").Value = "Cabildo"
$ws.Range("S268").Value = 2800
$ws.Range("T268").Value = 1
